# Apply the recorded value corrections for each sheet's H:N (pricing) columns.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).Value = ""
# Row 17
$ws.Cells.Item(17, 8).Value = 1251812.4
$ws.Cells.Item(17, 10).Value = 1668749.9
$ws.Cells.Item(17, 12).Value = 5006249.699999999
$ws.Cells.Item(17, 14).Value = -5006585.699999999
# Row 18
$ws.Cells.Item(18, 8).Value = 10950.1
$ws.Cells.Item(18, 9).Value = 10950.1
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 10950.1
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = -10666.1
$ws.Cells.Item(18, 14).Value = ""
# Row 31
$ws.Cells.Item(31, 8).Value = 337.83334
$ws.Cells.Item(31, 9).Value = 356.75
$ws.Cells.Item(31, 10).Value = 300
$ws.Cells.Item(31, 11).Value = 1070.25
$ws.Cells.Item(31, 12).Value = 900
$ws.Cells.Item(31, 13).Value = -840.25
$ws.Cells.Item(31, 14).Value = -1360
# Row 40
$ws.Cells.Item(40, 8).Value = 3756.25
$ws.Cells.Item(40, 10).Value = 5114.143
$ws.Cells.Item(40, 12).Value = 5114.143
$ws.Cells.Item(40, 14).Value = -5464.143
# Row 100
$ws.Cells.Item(100, 8).Value = 7038.0586
$ws.Cells.Item(100, 9).Value = 7895.364
$ws.Cells.Item(100, 10).Value = 5466.3335
$ws.Cells.Item(100, 11).Value = 7895.364
$ws.Cells.Item(100, 12).Value = 5466.3335
$ws.Cells.Item(100, 13).Value = -7354.364
$ws.Cells.Item(100, 14).Value = -6548.3335
# Row 112
$ws.Cells.Item(112, 8).Value = 1209.625
$ws.Cells.Item(112, 10).Value = 1206.2174
$ws.Cells.Item(112, 12).Value = 3618.6522
$ws.Cells.Item(112, 14).Value = -5834.6522
# Row 113
$ws.Cells.Item(113, 8).Value = 104899.3
$ws.Cells.Item(113, 10).Value = 5624.125
$ws.Cells.Item(113, 12).Value = 5624.125
$ws.Cells.Item(113, 14).Value = -12132.125
# Row 132
$ws.Cells.Item(132, 8).Value = 1579.3715
$ws.Cells.Item(132, 9).Value = 1579.3715
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4738.1145
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2208.1145
$ws.Cells.Item(132, 14).Value = ""
# Row 135
$ws.Cells.Item(135, 8).Value = 84963.164
$ws.Cells.Item(135, 9).Value = 1889
$ws.Cells.Item(135, 11).Value = 17001
$ws.Cells.Item(135, 13).Value = -14466
# Row 137
$ws.Cells.Item(137, 8).Value = 2029.3235
$ws.Cells.Item(137, 9).Value = 1840.5555
$ws.Cells.Item(137, 11).Value = 5521.666499999999
$ws.Cells.Item(137, 13).Value = -2971.666499999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1285.36
$ws.Cells.Item(2, 9).Value = 1159.1
$ws.Cells.Item(2, 11).Value = 1159.1
$ws.Cells.Item(2, 13).Value = -1046.1
# Row 32
$ws.Cells.Item(32, 8).Value = 4750.7827
$ws.Cells.Item(32, 9).Value = 4966.6665
$ws.Cells.Item(32, 10).Value = 2484
$ws.Cells.Item(32, 11).Value = 4966.6665
$ws.Cells.Item(32, 12).Value = 2484
$ws.Cells.Item(32, 13).Value = -4679.6665
$ws.Cells.Item(32, 14).Value = -3058
# Row 64
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).Value = ""
# Row 67
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).Value = ""
# Row 68
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = ""
# Row 71
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = ""
# Row 97
$ws.Cells.Item(97, 8).Value = 3500
$ws.Cells.Item(97, 9).Value = 3500
$ws.Cells.Item(97, 11).Value = 3500
$ws.Cells.Item(97, 13).Value = -3004
# Row 116
$ws.Cells.Item(116, 8).Value = 1285.36
$ws.Cells.Item(116, 9).Value = 1159.1
$ws.Cells.Item(116, 11).Value = 1159.1
$ws.Cells.Item(116, 13).Value = 1134.9
# Row 122
$ws.Cells.Item(122, 8).Value = 2387.0952
$ws.Cells.Item(122, 9).Value = 2390.4707
$ws.Cells.Item(122, 11).Value = 7171.4121
$ws.Cells.Item(122, 13).Value = -4721.4121

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1285.36
$ws.Cells.Item(3, 9).Value = 1159.1
$ws.Cells.Item(3, 11).Value = 1159.1
$ws.Cells.Item(3, 13).Value = -1045.1
# Row 86
$ws.Cells.Item(86, 8).Value = 13690.5
$ws.Cells.Item(86, 9).Value = 14285.75
$ws.Cells.Item(86, 11).Value = 14285.75
$ws.Cells.Item(86, 13).Value = -13162.75
# Row 89
$ws.Cells.Item(89, 8).Value = 13690.5
$ws.Cells.Item(89, 9).Value = 14285.75
$ws.Cells.Item(89, 11).Value = 71428.75
$ws.Cells.Item(89, 13).Value = -65812.75
# Row 134
$ws.Cells.Item(134, 8).Value = 1503.4667
$ws.Cells.Item(134, 9).Value = 1226.4445
$ws.Cells.Item(134, 11).Value = 3679.3335
$ws.Cells.Item(134, 13).Value = -1144.3335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Cells.Item(74, 8).Value = 59995
$ws.Cells.Item(74, 10).Value = 59995
$ws.Cells.Item(74, 12).Value = 59995
$ws.Cells.Item(74, 14).Value = -61743
# Row 77
$ws.Cells.Item(77, 8).Value = 59995
$ws.Cells.Item(77, 10).Value = 59995
$ws.Cells.Item(77, 12).Value = 179985
$ws.Cells.Item(77, 14).Value = -188721
# Row 97
$ws.Cells.Item(97, 8).Value = 30908.732
$ws.Cells.Item(97, 10).Value = 30908.732
$ws.Cells.Item(97, 12).Value = 30908.732
$ws.Cells.Item(97, 14).Value = -32890.732
# Row 132
$ws.Cells.Item(132, 8).Value = 4627.091
$ws.Cells.Item(132, 9).Value = 4322
$ws.Cells.Item(132, 10).Value = 6000
$ws.Cells.Item(132, 11).Value = 12966
$ws.Cells.Item(132, 12).Value = 18000
$ws.Cells.Item(132, 13).Value = -10436
$ws.Cells.Item(132, 14).Value = -23060
# Row 134
$ws.Cells.Item(134, 8).Value = 3198.9524
$ws.Cells.Item(134, 9).Value = 3030.0667
$ws.Cells.Item(134, 11).Value = 9090.2001
$ws.Cells.Item(134, 13).Value = -6555.2001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Cells.Item(6, 8).Value = 31
$ws.Cells.Item(6, 9).Value = 33
$ws.Cells.Item(6, 10).Value = 25
$ws.Cells.Item(6, 11).Value = 99
$ws.Cells.Item(6, 12).Value = 75
$ws.Cells.Item(6, 13).Value = 14
$ws.Cells.Item(6, 14).Value = -301
# Row 34
$ws.Cells.Item(34, 8).Value = 387.0909
$ws.Cells.Item(34, 9).Value = 164.75
$ws.Cells.Item(34, 10).Value = 980
$ws.Cells.Item(34, 11).Value = 494.25
$ws.Cells.Item(34, 12).Value = 2940
$ws.Cells.Item(34, 13).Value = -410.25
$ws.Cells.Item(34, 14).Value = -3108
# Row 86
$ws.Cells.Item(86, 8).Value = 583.75
$ws.Cells.Item(86, 9).Value = 465.1
$ws.Cells.Item(86, 11).Value = 1395.3
$ws.Cells.Item(86, 13).Value = -209.3000000000002
# Row 89
$ws.Cells.Item(89, 8).Value = 583.75
$ws.Cells.Item(89, 9).Value = 465.1
$ws.Cells.Item(89, 11).Value = 4185.900000000001
$ws.Cells.Item(89, 13).Value = 1742.099999999999
# Row 92
$ws.Cells.Item(92, 8).Value = 808.6667
$ws.Cells.Item(92, 9).Value = 880.4
$ws.Cells.Item(92, 10).Value = 450
$ws.Cells.Item(92, 11).Value = 2641.2
$ws.Cells.Item(92, 12).Value = 1350
$ws.Cells.Item(92, 13).Value = -1393.2
$ws.Cells.Item(92, 14).Value = -3846
# Row 100
$ws.Cells.Item(100, 8).Value = 8837
$ws.Cells.Item(100, 10).Value = 9567.799999999999
$ws.Cells.Item(100, 12).Value = 28703.4
$ws.Cells.Item(100, 14).Value = -30325.4
# Row 138
$ws.Cells.Item(138, 8).Value = 2427.5
$ws.Cells.Item(138, 9).Value = 1149.625
$ws.Cells.Item(138, 10).Value = 4983.25
$ws.Cells.Item(138, 11).Value = 3448.875
$ws.Cells.Item(138, 12).Value = 14949.75
$ws.Cells.Item(138, 13).Value = 1691.125
$ws.Cells.Item(138, 14).Value = -25229.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 2000
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 14).Value = -7900
# Row 126
$ws.Cells.Item(126, 8).Value = 4327.4287
$ws.Cells.Item(126, 9).Value = 3973.75
$ws.Cells.Item(126, 11).Value = 11921.25
$ws.Cells.Item(126, 13).Value = -9451.25
# Row 133
$ws.Cells.Item(133, 8).Value = 104935.25
$ws.Cells.Item(133, 10).Value = 104935.25
$ws.Cells.Item(133, 12).Value = 104935.25
$ws.Cells.Item(133, 14).Value = -115055.25
# Row 136
$ws.Cells.Item(136, 8).Value = 35191.867
$ws.Cells.Item(136, 10).Value = 35191.867
$ws.Cells.Item(136, 12).Value = 105575.601
$ws.Cells.Item(136, 14).Value = -110675.601

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 140
$ws.Cells.Item(140, 8).Value = 92428
$ws.Cells.Item(140, 10).Value = 92428
$ws.Cells.Item(140, 12).Value = 92428
$ws.Cells.Item(140, 14).Value = -102788

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Cells.Item(70, 8).Value = 20089.285
$ws.Cells.Item(70, 10).Value = 20089.285
$ws.Cells.Item(70, 12).Value = 20089.285
$ws.Cells.Item(70, 14).Value = -20719.285
# Row 73
$ws.Cells.Item(73, 8).Value = 20089.285
$ws.Cells.Item(73, 10).Value = 20089.285
$ws.Cells.Item(73, 12).Value = 20089.285
$ws.Cells.Item(73, 14).Value = -22273.285
